$wb = $excel.ActiveWorkbook

# --- Sheet 1: Biomass ---
$ws1 = $wb.Worksheets.Item("Biomass")
$ws1.Cells.Item(1,5).Value = "2020_Cod"
$ws1.Cells.Item(2,3).Value = 403588
$ws1.Cells.Item(2,5).Value = 364209
$ws1.Cells.Item(3,3).Value = 422439
$ws1.Cells.Item(3,5).Value = 405023
$ws1.Cells.Item(4,3).Value = 504136
$ws1.Cells.Item(4,5).Value = 471526
$ws1.Cells.Item(5,3).Value = 593197
$ws1.Cells.Item(5,5).Value = 530446
$ws1.Cells.Item(6,3).Value = 635060
$ws1.Cells.Item(6,5).Value = 556366
$ws1.Cells.Item(7,3).Value = 668967
$ws1.Cells.Item(7,5).Value = 579603
$ws1.Cells.Item(8,3).Value = 713828
$ws1.Cells.Item(8,5).Value = 618424
$ws1.Cells.Item(9,3).Value = 758519
$ws1.Cells.Item(9,5).Value = 651063
$ws1.Cells.Item(10,3).Value = 798787
$ws1.Cells.Item(10,5).Value = 690522
$ws1.Cells.Item(11,3).Value = 837433
$ws1.Cells.Item(11,5).Value = 742523
$ws1.Cells.Item(12,3).Value = 871227
$ws1.Cells.Item(12,5).Value = 787267
$ws1.Cells.Item(13,3).Value = 873994
$ws1.Cells.Item(13,5).Value = 794280
$ws1.Cells.Item(14,3).Value = 857974
$ws1.Cells.Item(14,5).Value = 787561
$ws1.Cells.Item(15,3).Value = 823846
$ws1.Cells.Item(15,5).Value = 769439
$ws1.Cells.Item(16,3).Value = 776061
$ws1.Cells.Item(16,5).Value = 736250
$ws1.Cells.Item(17,3).Value = 743411
$ws1.Cells.Item(17,5).Value = 707463
$ws1.Cells.Item(18,3).Value = 705219
$ws1.Cells.Item(18,5).Value = 665035
$ws1.Cells.Item(19,3).Value = 667998
$ws1.Cells.Item(19,5).Value = 628624
$ws1.Cells.Item(20,3).Value = 613486
$ws1.Cells.Item(20,5).Value = 585143
$ws1.Cells.Item(21,3).Value = 533819
$ws1.Cells.Item(21,5).Value = 519793
$ws1.Cells.Item(22,3).Value = 470311
$ws1.Cells.Item(22,5).Value = 465313
$ws1.Cells.Item(23,3).Value = 418822
$ws1.Cells.Item(23,5).Value = 416933
$ws1.Cells.Item(24,3).Value = 379079
$ws1.Cells.Item(24,5).Value = 377406
$ws1.Cells.Item(25,3).Value = 338390
$ws1.Cells.Item(25,5).Value = 338761
$ws1.Cells.Item(26,3).Value = 325728
$ws1.Cells.Item(26,5).Value = 324757
$ws1.Cells.Item(27,3).Value = 334660
$ws1.Cells.Item(27,5).Value = 328461
$ws1.Cells.Item(28,3).Value = 336484
$ws1.Cells.Item(28,5).Value = 323780
$ws1.Cells.Item(29,3).Value = 317059
$ws1.Cells.Item(29,5).Value = 300952
$ws1.Cells.Item(30,3).Value = 288103
$ws1.Cells.Item(30,5).Value = 274389
$ws1.Cells.Item(31,3).Value = 272454
$ws1.Cells.Item(31,5).Value = 266500
$ws1.Cells.Item(32,3).Value = 281250
$ws1.Cells.Item(32,5).Value = 284405
$ws1.Cells.Item(33,3).Value = 316237
$ws1.Cells.Item(33,5).Value = 321322
$ws1.Cells.Item(34,3).Value = 364318
$ws1.Cells.Item(34,5).Value = 365206
$ws1.Cells.Item(35,3).Value = 421953
$ws1.Cells.Item(35,5).Value = 412065
$ws1.Cells.Item(36,3).Value = 441055
$ws1.Cells.Item(36,5).Value = 425550
$ws1.Cells.Item(37,3).Value = 442457
$ws1.Cells.Item(37,5).Value = 431913
$ws1.Cells.Item(38,3).Value = 467633
$ws1.Cells.Item(38,5).Value = 475139
$ws1.Cells.Item(39,3).Value = 541959
$ws1.Cells.Item(39,5).Value = 553456
$ws1.Cells.Item(40,3).Value = 413621
$ws1.Cells.Item(40,5).Value = 408888
$ws1.Cells.Item(41,3).Value = 278457
$ws1.Cells.Item(41,5).Value = 266608
$ws1.Cells.Item(42,3).Value = 166636
$ws1.Cells.Item(42,5).Value = 157240
$ws1.Cells.Item(43,3).Value = 146433
$ws1.Cells.Item(43,5).Value = 131650
$ws1.Cells.Item(44,1).Value = 2019
$ws1.Cells.Item(44,5).Value = 149969
$ws1.Cells.Item(45,1).Value = 2020
$ws1.Cells.Item(45,5).Value = 186666

# --- Sheet 2: SSB ---
$ws2 = $wb.Worksheets.Item("SSB")
$ws2.Cells.Item(1,5).Value = "2020_Cod"
$ws2.Cells.Item(2,3).Value = 240906
$ws2.Cells.Item(2,5).Value = 220822
$ws2.Cells.Item(3,3).Value = 260534
$ws2.Cells.Item(3,5).Value = 239697
$ws2.Cells.Item(4,3).Value = 252020
$ws2.Cells.Item(4,5).Value = 235975
$ws2.Cells.Item(5,3).Value = 247466
$ws2.Cells.Item(5,5).Value = 240527
$ws2.Cells.Item(6,3).Value = 302871
$ws2.Cells.Item(6,5).Value = 287555
$ws2.Cells.Item(7,3).Value = 376994
$ws2.Cells.Item(7,5).Value = 333569
$ws2.Cells.Item(8,3).Value = 395472
$ws2.Cells.Item(8,5).Value = 344749
$ws2.Cells.Item(9,3).Value = 400666
$ws2.Cells.Item(9,5).Value = 347180
$ws2.Cells.Item(10,3).Value = 436257
$ws2.Cells.Item(10,5).Value = 377083
$ws2.Cells.Item(11,3).Value = 484999
$ws2.Cells.Item(11,5).Value = 422196
$ws2.Cells.Item(12,3).Value = 508411
$ws2.Cells.Item(12,5).Value = 454985
$ws2.Cells.Item(13,3).Value = 510659
$ws2.Cells.Item(13,5).Value = 466017
$ws2.Cells.Item(14,3).Value = 526359
$ws2.Cells.Item(14,5).Value = 486343
$ws2.Cells.Item(15,3).Value = 521522
$ws2.Cells.Item(15,5).Value = 486461
$ws2.Cells.Item(16,3).Value = 473885
$ws2.Cells.Item(16,5).Value = 446601
$ws2.Cells.Item(17,3).Value = 430265
$ws2.Cells.Item(17,5).Value = 410158
$ws2.Cells.Item(18,3).Value = 398097
$ws2.Cells.Item(18,5).Value = 382899
$ws2.Cells.Item(19,3).Value = 401249
$ws2.Cells.Item(19,5).Value = 388969
$ws2.Cells.Item(20,3).Value = 402597
$ws2.Cells.Item(20,5).Value = 388821
$ws2.Cells.Item(21,3).Value = 361454
$ws2.Cells.Item(21,5).Value = 345532
$ws2.Cells.Item(22,3).Value = 302929
$ws2.Cells.Item(22,5).Value = 293694
$ws2.Cells.Item(23,3).Value = 245753
$ws2.Cells.Item(23,5).Value = 247246
$ws2.Cells.Item(24,3).Value = 214551
$ws2.Cells.Item(24,5).Value = 220957
$ws2.Cells.Item(25,3).Value = 190885
$ws2.Cells.Item(25,5).Value = 194687
$ws2.Cells.Item(26,3).Value = 175240
$ws2.Cells.Item(26,5).Value = 175784
$ws2.Cells.Item(27,3).Value = 165709
$ws2.Cells.Item(27,5).Value = 167020
$ws2.Cells.Item(28,3).Value = 165570
$ws2.Cells.Item(28,5).Value = 165756
$ws2.Cells.Item(29,3).Value = 171104
$ws2.Cells.Item(29,5).Value = 166849
$ws2.Cells.Item(30,3).Value = 166219
$ws2.Cells.Item(30,5).Value = 158075
$ws2.Cells.Item(31,3).Value = 152138
$ws2.Cells.Item(31,5).Value = 141916
$ws2.Cells.Item(32,3).Value = 133144
$ws2.Cells.Item(32,5).Value = 124747
$ws2.Cells.Item(33,3).Value = 118934
$ws2.Cells.Item(33,5).Value = 116691
$ws2.Cells.Item(34,3).Value = 124955
$ws2.Cells.Item(34,5).Value = 126048
$ws2.Cells.Item(35,3).Value = 162165
$ws2.Cells.Item(35,5).Value = 164317
$ws2.Cells.Item(36,3).Value = 190668
$ws2.Cells.Item(36,5).Value = 186628
$ws2.Cells.Item(37,3).Value = 210816
$ws2.Cells.Item(37,5).Value = 198720
$ws2.Cells.Item(38,3).Value = 219493
$ws2.Cells.Item(38,5).Value = 205243
$ws2.Cells.Item(39,3).Value = 219628
$ws2.Cells.Item(39,5).Value = 213549
$ws2.Cells.Item(40,3).Value = 152560
$ws2.Cells.Item(40,5).Value = 156531
$ws2.Cells.Item(41,3).Value = 120170
$ws2.Cells.Item(41,5).Value = 125791
$ws2.Cells.Item(42,3).Value = 90748.9
$ws2.Cells.Item(42,5).Value = 89922.4
$ws2.Cells.Item(43,3).Value = 79446.3
$ws2.Cells.Item(43,5).Value = 71879.600000000006
$ws2.Cells.Item(44,1).Value = 2019
$ws2.Cells.Item(44,5).Value = 69588.3
$ws2.Cells.Item(45,1).Value = 2020
$ws2.Cells.Item(45,5).Value = 69262.5

# --- Sheet 3: R ---
$ws3 = $wb.Worksheets.Item("R")
$ws3.Cells.Item(1,5).Value = "2020_Cod"
$ws3.Cells.Item(2,5).Value = 1207580
$ws3.Cells.Item(3,5).Value = 377556
$ws3.Cells.Item(4,5).Value = 369733
$ws3.Cells.Item(5,5).Value = 624014
$ws3.Cells.Item(6,5).Value = 689951
$ws3.Cells.Item(7,5).Value = 756252
$ws3.Cells.Item(8,5).Value = 538912
$ws3.Cells.Item(9,5).Value = 709138
$ws3.Cells.Item(10,5).Value = 886695
$ws3.Cells.Item(11,5).Value = 499375
$ws3.Cells.Item(12,5).Value = 588083
$ws3.Cells.Item(13,5).Value = 597962
$ws3.Cells.Item(14,5).Value = 632229
$ws3.Cells.Item(15,5).Value = 749185
$ws3.Cells.Item(16,5).Value = 444758
$ws3.Cells.Item(17,5).Value = 385255
$ws3.Cells.Item(18,5).Value = 309854
$ws3.Cells.Item(19,5).Value = 347856
$ws3.Cells.Item(20,5).Value = 438067
$ws3.Cells.Item(21,5).Value = 309470
$ws3.Cells.Item(22,5).Value = 293505
$ws3.Cells.Item(23,5).Value = 272155
$ws3.Cells.Item(24,5).Value = 366527
$ws3.Cells.Item(25,5).Value = 439377
$ws3.Cells.Item(26,5).Value = 250745
$ws3.Cells.Item(27,5).Value = 193147
$ws3.Cells.Item(28,5).Value = 244348
$ws3.Cells.Item(29,5).Value = 307845
$ws3.Cells.Item(30,5).Value = 420358
$ws3.Cells.Item(31,5).Value = 686754
$ws3.Cells.Item(32,5).Value = 443195
$ws3.Cells.Item(33,5).Value = 651882
$ws3.Cells.Item(34,5).Value = 391813
$ws3.Cells.Item(35,5).Value = 506839
$ws3.Cells.Item(36,5).Value = 655108
$ws3.Cells.Item(37,5).Value = 1215110
$ws3.Cells.Item(38,5).Value = 638080
$ws3.Cells.Item(39,5).Value = 211074
$ws3.Cells.Item(40,5).Value = 260163
$ws3.Cells.Item(41,5).Value = 168038
$ws3.Cells.Item(42,5).Value = 246044
$ws3.Cells.Item(43,5).Value = 389895
$ws3.Cells.Item(44,1).Value = 2019
$ws3.Cells.Item(44,5).Value = 399011
$ws3.Cells.Item(45,1).Value = 2020
$ws3.Cells.Item(45,5).Value = 463705

